$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mark 3"
$ws.Range("A2").Value = "Mark 4"

$ws.Range("B6").Select()
